# UKIM ("People with Significant Control") template header wording refresh.
# Shorten/simplify the column-header labels in row 1 of Sheet1 and make the
# formatting of the header row consistent across its full width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (columns stay in the same places, wording changes) ---
$ws.Range("B1").Value = "Residential address"
$ws.Range("D1").Value = "National Insurance number"
$ws.Range("E1").Value = "Identification number if no National Insurance number (eg passport number, driver's licence, national identity card)"

# --- Make the rest of row 1 (F1:Z1) share the exact same formatting as the
#     populated header cells (A1:E1), instead of a slightly different bold
#     style, by copying the format from A1 across. ---
$ws.Range("A1").Copy()
$ws.Range("F1:Z1").PasteSpecial(-4122)
